$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Rename the worksheet (sheet tab name in workbook.xml)
$ws.Name = "alpha2F"

# Update the precise numeric values in row 13
$ws.Range("C13").Value = 0.9875976918588665
$ws.Range("F13").Value = 0.9875976918588665
$ws.Range("H13").Value = 0.9916696337485124
$ws.Range("L13").Value = 0.9909575396576098
